# "Add files via upload" - refresh the Custom Batting Order DAX measure on the
# Calculated_Columns sheet: the old formula referenced dim_player[name] (the
# playing-XI shortlist); it is replaced with a SWITCH() keyed on
# dim_players_wc2023[name] for the actual finalists' batting order. The
# "Table" cell for that row is retouched too, and the backing Excel Table
# (Table2) is resized one row taller to match the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calculated_Columns")
$ws.Activate()

$newFormula = @"
SWITCH(
TRUE(),
dim_players_wc2023[name] = "Rohit Sharma",1,
dim_players_wc2023[name] = "Quinton de Kock",2,
dim_players_wc2023[name] = "Virat Kohli",3,
dim_players_wc2023[name] = "Daryl Mitchell",4,
dim_players_wc2023[name] = "Rachin Ravindra" ,4,
dim_players_wc2023[name] = "KL Rahul" ,5,
dim_players_wc2023[name] = "Glenn Maxwell" ,6,
dim_players_wc2023[name] = "Glenn Phillips" ,7,
dim_players_wc2023[name] = "Ravindra Jadeja" ,8,
dim_players_wc2023[name] = "Mohammed Shami" ,9,
dim_players_wc2023[name] = "Jasprit Bumrah" ,10,
dim_players_wc2023[name] = "Adam Zampa" ,11,
dim_players_wc2023[name] = "Gerald Coetzee" ,11,
)
"@

# The "Table" column for that same row picked up a stray retype
$ws.Range("E4").Value = "dim_players_wc2023s_wc2023"

# Row 4 = Sno 3 = "Custom Batting Order"
$ws.Range("D4").Value = $newFormula
$ws.Range("D4").WrapText = $true
$ws.Range("D4").Select()
$ws.Range("D4").EntireRow.RowHeight = 230.4

# Grow the Table2 list object by one (blank) row, as happened in the upload
$lo = $ws.ListObjects.Item("Table2")
$lo.Resize($ws.Range("A1:E5"))

$ws.Range("D4").Select()
